$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.940.97"
$ws.Range("E2").Value = "  -0.60%  "
Set-TextValue $ws.Range("D3") "1.952.99"
$ws.Range("E3").Value = "  -0.79%  "
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  -0.46%  "
Set-TextValue $ws.Range("D5") "242.54"
$ws.Range("E5").Value = "  -2.28%  "
Set-TextValue $ws.Range("D6") "1.000"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  +0.05%  "
Set-TextValue $ws.Range("D8") "0.2936"
$ws.Range("E8").Value = "  -0.73%  "
Set-TextValue $ws.Range("D9") "0.06949"
$ws.Range("E9").Value = "  +1.71%  "
Set-TextValue $ws.Range("D10") "19.47"
$ws.Range("E10").Value = "  +1.36%  "
Set-TextValue $ws.Range("D11") "106.96"
$ws.Range("E11").Value = "  -0.58%  "
Set-TextValue $ws.Range("D12") "1.948.78"
$ws.Range("E12").Value = "  -0.99%  "
Set-TextValue $ws.Range("D13") "0.07757"
$ws.Range("E13").Value = "  -0.45%  "
Set-TextValue $ws.Range("D14") "5.342"
$ws.Range("E14").Value = "  -1.85%  "
Set-TextValue $ws.Range("D15") "0.6952"
$ws.Range("E15").Value = "  -1.66%  "
Set-TextValue $ws.Range("D16") "280.22"
$ws.Range("E16").Value = "  -1.81%  "
Set-TextValue $ws.Range("D17") "30.944.12"
$ws.Range("E17").Value = "  -0.62%  "
Set-TextValue $ws.Range("D18") "0.000007756"
$ws.Range("E18").Value = "  +0.23%  "
Set-TextValue $ws.Range("D19") "13.20"
$ws.Range("E19").Value = "  -0.62%  "
Set-TextValue $ws.Range("D20") "2.231.20"
$ws.Range("E20").Value = "  +0.37%  "
Set-TextValue $ws.Range("D21") "1.001"
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("E22").Value = "  -2.56%  "
Set-TextValue $ws.Range("D23") "1.000"
$ws.Range("E23").Value = "  -0.54%  "
Set-TextValue $ws.Range("D24") "6.489"
$ws.Range("E24").Value = "  -2.17%  "
Set-TextValue $ws.Range("D25") "9.731"
$ws.Range("E25").Value = "  -2.82%  "
Set-TextValue $ws.Range("D26") "168.13"
$ws.Range("E26").Value = "  -1.20%  "
Set-TextValue $ws.Range("D27") "19.67"
$ws.Range("E27").Value = "  -2.14%  "
Set-TextValue $ws.Range("D28") "2.169"
$ws.Range("E28").Value = "  -0.97%  "
Set-TextValue $ws.Range("D29") "0.1041"
$ws.Range("E29").Value = "  -2.36%  "
Set-TextValue $ws.Range("D30") "1.395"
$ws.Range("E30").Value = "  -3.62%  "
Set-TextValue $ws.Range("D31") "4.583"
$ws.Range("E31").Value = "  -5.37%  "
Set-TextValue $ws.Range("D32") "1.553"
$ws.Range("E32").Value = "  -2.95%  "
Set-TextValue $ws.Range("D33") "4.395"
$ws.Range("E33").Value = "  -2.51%  "
Set-TextValue $ws.Range("D34") "0.04869"
$ws.Range("E34").Value = "  -4.37%  "
Set-TextValue $ws.Range("D35") "0.7507"
$ws.Range("E35").Value = "  -2.63%  "
Set-TextValue $ws.Range("D36") "1.163"
$ws.Range("E36").Value = "  -0.60%  "
Set-TextValue $ws.Range("D37") "2.725"
$ws.Range("E37").Value = "  -0.48%  "
Set-TextValue $ws.Range("D38") "0.01995"
$ws.Range("E38").Value = "  -2.32%  "
Set-TextValue $ws.Range("D39") "2.679"
$ws.Range("E39").Value = "  -2.16%  "
Set-TextValue $ws.Range("D40") "6.505"
$ws.Range("E40").Value = "  +0.62%  "
Set-TextValue $ws.Range("D41") "77.71"
$ws.Range("E41").Value = "  +6.11%  "
Set-TextValue $ws.Range("D42") "2.107"
$ws.Range("E42").Value = "  -1.15%  "
Set-TextValue $ws.Range("D43") "0.8953"
$ws.Range("E43").Value = "  +0.96%  "
Set-TextValue $ws.Range("D44") "108.76"
$ws.Range("E44").Value = "  -1.17%  "
Set-TextValue $ws.Range("D45") "0.4430"
$ws.Range("E45").Value = "  -0.87%  "
Set-TextValue $ws.Range("D46") "1.0000"
$ws.Range("E46").Value = "  -0.34%  "
Set-TextValue $ws.Range("D47") "7.746"
$ws.Range("E47").Value = "  +3.00%  "
Set-TextValue $ws.Range("D48") "993.68"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("E49").Value = "  -1.98%  "
Set-TextValue $ws.Range("D50") "9.240"
$ws.Range("E50").Value = "  -1.83%  "
Set-TextValue $ws.Range("D51") "35.82"
$ws.Range("E51").Value = "  -0.53%  "
